$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Real")
$ws2 = $wb.Worksheets.Item("Fake")

# Fill in new cell values in the exact order they were entered,
# so that shared-string indices come out in the same order as the target.
$ws2.Range("A13").Value = "https://100percentfedup.com/oops-traitor-and-former-trump-education-sec-betsy-devos-trashes-trump-while-simultaneously-funding-gop-governor-candidate-gunning-for-trump-endorsement/"
$ws2.Range("A14").Value = "https://100percentfedup.com/the-biden-regime-claimed-that-their-disinformation-board-was-just-for-advisory-purposes-newly-released-documents-paint-a-darker-picture/"
$ws2.Range("A15").Value = "https://100percentfedup.com/lol-democrat-darling-rep-liz-cheney-gets-roasted-on-social-media-during-jan-6-witch-hunt-show-trial-liz-is-lying/"
$ws2.Range("A16").Value = "https://100percentfedup.com/finally-fauci-hints-at-when-he-is-going-to-retire-from-government/"
$ws2.Range("A17").Value = "https://100percentfedup.com/transgender-former-amazon-employee-goes-to-trial-for-stealing-over-100-million-capitol-one-customers-information/"
$ws2.Range("A18").Value = "https://100percentfedup.com/new-documents-show-biden-regimes-plan-to-send-migrants-to-cities-deeper-inside-u-s/"
$ws2.Range("A19").Value = "https://100percentfedup.com/washington-post-column-admits-dems-will-lose-hard-in-midterms-but-heres-the-twist/"
$ws2.Range("A23").Value = "https://pieceofmindful.com/2020/04/06/bombshell-who-coronavirus-pcr-test-primer-sequence-is-found-in-all-human-dna/"
$ws2.Range("A20").Value = "https://21stcenturywire.com/2022/06/10/vernon-coleman-the-death-of-health-care-in-britain/"
$ws2.Range("A21").Value = "https://21stcenturywire.com/2022/06/09/disinformation-board-leaked-documents-expose-agenda-behind-dhs-ministry-of-truth/"
$ws2.Range("A22").Value = "https://21stcenturywire.com/2022/06/08/was-the-pandemic-the-perfect-cover-for-the-great-reset/"
$ws2.Range("A24").Value = "https://21stcenturywire.com/2022/05/20/biden-orders-millions-of-monkeypox-vaccines-after-one-man-allegedly-infected/"
$ws2.Range("A25").Value = "https://21stcenturywire.com/2022/05/21/vernon-coleman-the-clock-is-ticking-loudly/"
$ws2.Range("A26").Value = "https://www.activistpost.com/2022/06/pfizer-ceo-head-of-cia-facebook-vp-other-elites-secretly-meeting-in-dc-corporate-media-is-silent.html"
$ws2.Range("A27").Value = "https://americanlookout.com/report-new-management-at-cnn-preparing-to-fire-partisan-hosts/"
$ws2.Range("A28").Value = "https://americanlookout.com/propaganda-democrats-hire-professional-tv-producer-for-prime-time-january-6th-hearings/"
$ws2.Range("A29").Value = "https://americanlookout.com/poll-finds-most-americans-believe-trump-is-not-responsible-for-january-6th/"
$ws2.Range("A30").Value = "https://americanlookout.com/report-john-deere-moving-part-of-production-from-iowa-to-mexico/"
$ws1.Range("A11").Value = "https://edition.cnn.com/2022/05/26/business/russia-economy-ruble-interest-rates/index.html"
$ws1.Range("A13").Value = "https://www.cnbc.com/2022/06/09/stock-market-news-open-to-close.html"
$ws1.Range("A14").Value = "https://www.cnbc.com/2022/06/10/flu-hepatitis-monkeypox-diseases-suppressed-during-covid-are-back.html"
$ws1.Range("A15").Value = "https://www.cnbc.com/2022/06/10/investment-banks-say-its-time-to-get-back-into-china-with-goldman-naming-10-top-stocks.html"
$ws1.Range("A16").Value = "https://www.cnbc.com/2022/06/10/klarna-ceo-defends-business-despite-massive-losses-and-layoffs.html"
$ws1.Range("A12").Value = "https://edition.cnn.com/2022/06/09/americas/dom-phillips-bruno-pereira-missing-suspect-intl-latam/index.html"
$ws1.Range("A17").Value = "https://www.cnbc.com/2022/06/10/ai-gurus-are-leaving-big-tech-to-work-on-buzzy-new-start-ups.html"
$ws1.Range("A18").Value = "https://www.cnbc.com/2022/06/10/carrie-lam-says-hong-kong-hasnt-become-just-another-chinese-city.html"
$ws1.Range("A21").Value = "https://edition.cnn.com/2022/06/08/americas/maduro-ankara-analysis-latam/index.html"
$ws1.Range("A19").Value = "https://www.cnbc.com/2022/06/10/tesla-cancels-three-june-online-hiring-events-for-china.html"
$ws1.Range("A22").Value = "https://edition.cnn.com/2022/06/10/politics/alejandro-mayorkas-interview-cnntv/index.html"
$ws1.Range("A20").Value = "https://www.politico.com/news/2022/06/09/cheney-scott-perry-jan-6-hearing-00038724"
$ws1.Range("A23").Value = "https://www.washingtonpost.com/world/2022/06/10/russia-putin-peter-the-great-ukraine-war/?itid=mr_world_2"
$ws1.Range("A24").Value = "https://www.washingtonpost.com/world/2022/06/09/china-rumors-xi-covid-politics/?itid=mr_world_3"
$ws1.Range("A25").Value = "https://www.washingtonpost.com/world/2022/06/09/who-sago-covid-origins/?itid=mr_world_4"
$ws1.Range("A27").Value = "https://www.politico.com/news/magazine/2022/06/10/ruth-bader-ginsburg-retire-legacy-00038638"
$ws1.Range("A26").Value = "https://www.washingtonpost.com/world/2022/06/09/belgium-king-philippe-democratic-republic-congo-colonialism/?itid=mr_world_5"
$ws1.Range("A28").Value = "https://www.politico.com/news/2022/06/10/chesa-boudin-progressive-agenda-california-00038675"
$ws1.Range("A29").Value = "https://www.politico.com/news/2022/06/09/biden-abortion-00038565"
$ws1.Range("A30").Value = "https://www.politico.com/sponsored-content/2022/05/seizing-the-next-revolution?utm_source=native&utm_medium=hp"

# Add hyperlinks in the same order as the original edit, so relationship ids line up.
$ws1.Hyperlinks.Add($ws1.Range("A11"), "https://edition.cnn.com/2022/05/26/business/russia-economy-ruble-interest-rates/index.html") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A13"), "https://www.cnbc.com/2022/06/09/stock-market-news-open-to-close.html") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A14"), "https://www.cnbc.com/2022/06/10/flu-hepatitis-monkeypox-diseases-suppressed-during-covid-are-back.html") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A15"), "https://www.cnbc.com/2022/06/10/investment-banks-say-its-time-to-get-back-into-china-with-goldman-naming-10-top-stocks.html") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A16"), "https://www.cnbc.com/2022/06/10/klarna-ceo-defends-business-despite-massive-losses-and-layoffs.html") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A17"), "https://www.cnbc.com/2022/06/10/ai-gurus-are-leaving-big-tech-to-work-on-buzzy-new-start-ups.html") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A18"), "https://www.cnbc.com/2022/06/10/carrie-lam-says-hong-kong-hasnt-become-just-another-chinese-city.html") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A19"), "https://www.cnbc.com/2022/06/10/tesla-cancels-three-june-online-hiring-events-for-china.html") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A12"), "https://edition.cnn.com/2022/06/09/americas/dom-phillips-bruno-pereira-missing-suspect-intl-latam/index.html") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A22"), "https://edition.cnn.com/2022/06/10/politics/alejandro-mayorkas-interview-cnntv/index.html") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A21"), "https://edition.cnn.com/2022/06/08/americas/maduro-ankara-analysis-latam/index.html") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A20"), "https://www.politico.com/news/2022/06/09/cheney-scott-perry-jan-6-hearing-00038724") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A27"), "https://www.politico.com/news/magazine/2022/06/10/ruth-bader-ginsburg-retire-legacy-00038638") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A28"), "https://www.politico.com/news/2022/06/10/chesa-boudin-progressive-agenda-california-00038675") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A23"), "https://www.washingtonpost.com/world/2022/06/10/russia-putin-peter-the-great-ukraine-war/?itid=mr_world_2") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A29"), "https://www.politico.com/news/2022/06/09/biden-abortion-00038565") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A24"), "https://www.washingtonpost.com/world/2022/06/09/china-rumors-xi-covid-politics/?itid=mr_world_3") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A25"), "https://www.washingtonpost.com/world/2022/06/09/who-sago-covid-origins/?itid=mr_world_4") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A26"), "https://www.washingtonpost.com/world/2022/06/09/belgium-king-philippe-democratic-republic-congo-colonialism/?itid=mr_world_5") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A30"), "https://www.politico.com/sponsored-content/2022/05/seizing-the-next-revolution?utm_source=native&utm_medium=hp") | Out-Null

$ws2.Hyperlinks.Add($ws2.Range("A13"), "https://100percentfedup.com/oops-traitor-and-former-trump-education-sec-betsy-devos-trashes-trump-while-simultaneously-funding-gop-governor-candidate-gunning-for-trump-endorsement/") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A14"), "https://100percentfedup.com/the-biden-regime-claimed-that-their-disinformation-board-was-just-for-advisory-purposes-newly-released-documents-paint-a-darker-picture/") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A15"), "https://100percentfedup.com/lol-democrat-darling-rep-liz-cheney-gets-roasted-on-social-media-during-jan-6-witch-hunt-show-trial-liz-is-lying/") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A16"), "https://100percentfedup.com/finally-fauci-hints-at-when-he-is-going-to-retire-from-government/") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A17"), "https://100percentfedup.com/transgender-former-amazon-employee-goes-to-trial-for-stealing-over-100-million-capitol-one-customers-information/") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A18"), "https://100percentfedup.com/new-documents-show-biden-regimes-plan-to-send-migrants-to-cities-deeper-inside-u-s/") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A19"), "https://100percentfedup.com/washington-post-column-admits-dems-will-lose-hard-in-midterms-but-heres-the-twist/") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A23"), "https://pieceofmindful.com/2020/04/06/bombshell-who-coronavirus-pcr-test-primer-sequence-is-found-in-all-human-dna/") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A20"), "https://21stcenturywire.com/2022/06/10/vernon-coleman-the-death-of-health-care-in-britain/") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A21"), "https://21stcenturywire.com/2022/06/09/disinformation-board-leaked-documents-expose-agenda-behind-dhs-ministry-of-truth/") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A22"), "https://21stcenturywire.com/2022/06/08/was-the-pandemic-the-perfect-cover-for-the-great-reset/") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A24"), "https://21stcenturywire.com/2022/05/20/biden-orders-millions-of-monkeypox-vaccines-after-one-man-allegedly-infected/") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A25"), "https://21stcenturywire.com/2022/05/21/vernon-coleman-the-clock-is-ticking-loudly/") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A26"), "https://www.activistpost.com/2022/06/pfizer-ceo-head-of-cia-facebook-vp-other-elites-secretly-meeting-in-dc-corporate-media-is-silent.html") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A27"), "https://americanlookout.com/report-new-management-at-cnn-preparing-to-fire-partisan-hosts/") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A28"), "https://americanlookout.com/propaganda-democrats-hire-professional-tv-producer-for-prime-time-january-6th-hearings/") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A29"), "https://americanlookout.com/poll-finds-most-americans-believe-trump-is-not-responsible-for-january-6th/") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A30"), "https://americanlookout.com/report-john-deere-moving-part-of-production-from-iowa-to-mexico/") | Out-Null

# Restore the plain 'Hyperlink' cell style (Hyperlinks.Add() leaves a duplicated style behind).
$ws1.Range("A11:A30").Style = "Hyperlink"
$ws2.Range("A13:A30").Style = "Hyperlink"

# Update window scroll position / selection to match the saved view state.
$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws2.Range("A31").Select()

$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("A31").Select()
